$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric/percentage-looking text cells (columns D and E) must be forced to
# text format first so Excel does not auto-convert them to numbers and lose
# formatting (trailing zeros, percent signs, etc.), matching the original
# inline-string cell contents.
$deCells = @("D2","E2","E3","E4","D5","E5","E6","D7","E7","D8","E8","D9","E9","D10","E10","D11","E11","D12","E12","D13","E13","D14","E14","D15","E15","D16","E16","D17","E17","D18","E18","D19","E19","D20","E20","D21","E21","D22","E22","D23","E23","D24","E24","D25","E25","D26","E26","D27","E27","E28","D40","E40","D41","E41","D42","E42","D43","E43","D44","E44","D45","E45","D46","E46","E47")
foreach ($addr in $deCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated price / volume values
$ws.Range("D2").Value = "285.96"
$ws.Range("E2").Value = "0.97%"
$ws.Range("E3").Value = "3.02%"
$ws.Range("E4").Value = "0.83%"
$ws.Range("D5").Value = "0.06721"
$ws.Range("E5").Value = "3.10%"
$ws.Range("E6").Value = "0.94%"
$ws.Range("D7").Value = "1.395"
$ws.Range("E7").Value = "-1.54%"
$ws.Range("D8").Value = "0.8989"
$ws.Range("E8").Value = "-2.60%"
$ws.Range("D9").Value = "0.1576"
$ws.Range("E9").Value = "1.91%"
$ws.Range("D10").Value = "0.06886"
$ws.Range("E10").Value = "6.21%"
$ws.Range("D11").Value = "0.07626"
$ws.Range("E11").Value = "0.27%"
$ws.Range("D12").Value = "0.02919"
$ws.Range("E12").Value = "0.72%"
$ws.Range("D13").Value = "0.08989"
$ws.Range("E13").Value = "0.51%"
$ws.Range("D14").Value = "0.001585"
$ws.Range("E14").Value = "-0.16%"
$ws.Range("D15").Value = "0.04487"
$ws.Range("E15").Value = "1.61%"
$ws.Range("D16").Value = "0.0006474"
$ws.Range("E16").Value = "1.09%"
$ws.Range("D17").Value = "0.006472"
$ws.Range("E17").Value = "7.16%"
$ws.Range("D18").Value = "3.448"
$ws.Range("E18").Value = "0.00%"
$ws.Range("D19").Value = "3.434"
$ws.Range("E19").Value = "1.37%"
$ws.Range("D20").Value = "2.231"
$ws.Range("E20").Value = "-0.70%"
$ws.Range("D21").Value = "0.3205"
$ws.Range("E21").Value = "0.42%"
$ws.Range("D22").Value = "0.1319"
$ws.Range("E22").Value = "2.32%"
$ws.Range("D23").Value = "4.039"
$ws.Range("E23").Value = "1.19%"
$ws.Range("D24").Value = "0.1580"
$ws.Range("E24").Value = "2.54%"
$ws.Range("D25").Value = "0.001201"
$ws.Range("E25").Value = "1.10%"
$ws.Range("D26").Value = "0.004376"
$ws.Range("E26").Value = "-0.98%"
$ws.Range("D27").Value = "0.0001168"
$ws.Range("E27").Value = "-6.65%"
$ws.Range("E28").Value = "-0.07%"
$ws.Range("D40").Value = "0.04255"
$ws.Range("E40").Value = "2.40%"
$ws.Range("D41").Value = "0.006814"
$ws.Range("E41").Value = "2.12%"
$ws.Range("D42").Value = "0.1238"
$ws.Range("E42").Value = "0.99%"
$ws.Range("D43").Value = "0.002167"
$ws.Range("E43").Value = "-1.18%"
$ws.Range("D44").Value = "0.01152"
$ws.Range("E44").Value = "-5.19%"
$ws.Range("D45").Value = "0.00005730"
$ws.Range("E45").Value = "1.92%"
$ws.Range("D46").Value = "1.918"
$ws.Range("E46").Value = "-2.41%"
$ws.Range("E47").Value = "15.33%"

# Apply the updated coin name / link text values (rows shifted due to a new
# "CoinExToken" entry being inserted ahead of the existing rows)
$ws.Range("B15").Value = "CoinExToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("B16").Value = "One"
$ws.Range("C16").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("B19").Value = "GateToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("B20").Value = "BTSEToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("B21").Value = "BitpandaEcosystemToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("B22").Value = "ProBitToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("B23").Value = "MCDex"
$ws.Range("C23").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("B24").Value = "ZBToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
